$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CASA-Student List")

# Row 18 was the blank row just under the last filled student (row 17) inside
# the table's visual area. Fill it in with a new student record: Kenia Ramos.
$ws.Range("C18").Value = 14
$ws.Range("D18").Value = "Kenia Ramos"

# Email column gets the address as visible text plus a mailto: hyperlink.
$ws.Range("E18").Value = "Kramos@wearecasa.org"
$ws.Hyperlinks.Add($ws.Range("E18"), "mailto:Kramos@wearecasa.org")

# The remaining columns repeat the same placeholder labels used on every
# other data row (copied from row 17, the row directly above).
$ws.Range("F18").Value = $ws.Range("F17").Value2
$ws.Range("G18").Value = $ws.Range("G17").Value2
$ws.Range("H18").Value = $ws.Range("H17").Value2
$ws.Range("I18").Value = $ws.Range("I17").Value2
$ws.Range("J18").Value = $ws.Range("J17").Value2
$ws.Range("K18").Value = $ws.Range("K17").Value2
$ws.Range("L18").Value = $ws.Range("L17").Value2
